# Apply updated Chocobo leve profit calculations across sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 39899.5
$ws.Range("J3").Value = 39899.5
$ws.Range("L3").Value = 39899.5
$ws.Range("N3").Value = -40127.5

# Row 102
$ws.Range("H102").Value = 39899.5
$ws.Range("J102").Value = 39899.5
$ws.Range("L102").Value = 39899.5
$ws.Range("N102").Value = -46389.5

# Row 106
$ws.Range("H106").Value = 2249.75
$ws.Range("I106").Value = 1499.5
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1499.5
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -868.5
$ws.Range("N106").Value = -4262

# Row 141
$ws.Range("H141").Value = 277750.5
$ws.Range("I141").Value = 368667.34
$ws.Range("K141").Value = 1106002.02
$ws.Range("M141").Value = -1100822.02

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 8526.462
$ws.Range("I74").Value = 9761.111000000001
$ws.Range("J74").Value = 5748.5
$ws.Range("K74").Value = 9761.111000000001
$ws.Range("L74").Value = 5748.5
$ws.Range("M74").Value = -8887.111000000001
$ws.Range("N74").Value = -7496.5

# Row 76
$ws.Range("H76").Value = 29172.8
$ws.Range("J76").Value = 29172.8
$ws.Range("L76").Value = 29172.8
$ws.Range("N76").Value = -29848.8

# Row 77
$ws.Range("H77").Value = 8526.462
$ws.Range("I77").Value = 9761.111000000001
$ws.Range("J77").Value = 5748.5
$ws.Range("K77").Value = 48805.55500000001
$ws.Range("L77").Value = 28742.5
$ws.Range("M77").Value = -44437.55500000001
$ws.Range("N77").Value = -37478.5

# Row 79
$ws.Range("H79").Value = 29172.8
$ws.Range("J79").Value = 29172.8
$ws.Range("L79").Value = 29172.8
$ws.Range("N79").Value = -31512.8

# Row 103
$ws.Range("H103").Value = 35041.08
$ws.Range("J103").Value = 35041.08
$ws.Range("L103").Value = 35041.08
$ws.Range("N103").Value = -37385.08

# Row 112
$ws.Range("H112").Value = 30476.191
$ws.Range("J112").Value = 30476.191
$ws.Range("L112").Value = 30476.191
$ws.Range("N112").Value = -33430.191

$ws = $wb.Worksheets.Item("BSM")
# Row 69
$ws.Range("H69").Value = 29000
$ws.Range("J69").Value = 29000
$ws.Range("L69").Value = 29000
$ws.Range("N69").Value = -30622

# Row 72
$ws.Range("H72").Value = 29000
$ws.Range("J72").Value = 29000
$ws.Range("L72").Value = 87000
$ws.Range("N72").Value = -95112

$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 48600
$ws.Range("J52").Value = 48600
$ws.Range("L52").Value = 48600
$ws.Range("N52").Value = -49188

# Row 57
$ws.Range("H57").Value = 48011.8
$ws.Range("J57").Value = 48011.8
$ws.Range("L57").Value = 48011.8
$ws.Range("N57").Value = -49131.8

# Row 63
$ws.Range("H63").Value = 49995
$ws.Range("J63").Value = 49995
$ws.Range("L63").Value = 49995
$ws.Range("N63").Value = -51367

# Row 66
$ws.Range("H66").Value = 49995
$ws.Range("J66").Value = 49995
$ws.Range("L66").Value = 149985
$ws.Range("N66").Value = -156849

# Row 68
$ws.Range("H68").Value = 47676.1
$ws.Range("J68").Value = 47676.1
$ws.Range("L68").Value = 47676.1
$ws.Range("N68").Value = -49174.1

# Row 71
$ws.Range("H71").Value = 47676.1
$ws.Range("J71").Value = 47676.1
$ws.Range("L71").Value = 143028.3
$ws.Range("N71").Value = -150516.3

# Row 109
$ws.Range("H109").Value = 34799.25
$ws.Range("J109").Value = 34799.25
$ws.Range("L109").Value = 34799.25
$ws.Range("N109").Value = -36879.25

# Row 137
$ws.Range("H137").Value = 41922.5
$ws.Range("J137").Value = 41922.5
$ws.Range("L137").Value = 41922.5
$ws.Range("N137").Value = -52122.5

$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -15256

# Row 108
$ws.Range("H108").Value = 2000
$ws.Range("I108").Value = 2000
$ws.Range("K108").Value = 6000
$ws.Range("M108").Value = -3120

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 36474.875
$ws.Range("J15").Value = 36474.875
$ws.Range("L15").Value = 36474.875
$ws.Range("N15").Value = -37050.875

# Row 41
$ws.Range("H41").Value = 9958.066000000001
$ws.Range("I41").Value = 1437.625
$ws.Range("J41").Value = 19695.715
$ws.Range("K41").Value = 1437.625
$ws.Range("L41").Value = 19695.715
$ws.Range("M41").Value = -1082.625
$ws.Range("N41").Value = -20405.715

# Row 81
$ws.Range("H81").Value = 36474.875
$ws.Range("J81").Value = 36474.875
$ws.Range("L81").Value = 36474.875
$ws.Range("N81").Value = -38470.875

# Row 84
$ws.Range("H84").Value = 36474.875
$ws.Range("J84").Value = 36474.875
$ws.Range("L84").Value = 109424.625
$ws.Range("N84").Value = -119408.625

# Row 107
$ws.Range("H107").Value = 8547664
$ws.Range("I107").Value = 256.8
$ws.Range("J107").Value = 13889794
$ws.Range("K107").Value = 256.8
$ws.Range("L107").Value = 13889794
$ws.Range("M107").Value = 1663.2
$ws.Range("N107").Value = -13893634

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 951.6
$ws.Range("J2").Value = 1002
$ws.Range("L2").Value = 1002
$ws.Range("N2").Value = -1226

# Row 110
$ws.Range("H110").Value = 39800
$ws.Range("J110").Value = 39800
$ws.Range("L110").Value = 39800
$ws.Range("N110").Value = -47980

# Row 127
$ws.Range("H127").Value = 26704.666
$ws.Range("J127").Value = 26704.666
$ws.Range("L127").Value = 26704.666
$ws.Range("N127").Value = -36624.666

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 72856.55499999999
$ws.Range("J46").Value = 72856.55499999999
$ws.Range("L46").Value = 72856.55499999999
$ws.Range("N46").Value = -73318.55499999999

# Row 57
$ws.Range("H57").Value = 17950.2
$ws.Range("J57").Value = 17950.2
$ws.Range("L57").Value = 17950.2
$ws.Range("N57").Value = -19458.2

# Row 101
$ws.Range("H101").Value = 29799.5
$ws.Range("J101").Value = 29799.5
$ws.Range("L101").Value = 29799.5
$ws.Range("N101").Value = -36289.5

# Row 122
$ws.Range("H122").Value = 2307.7334
$ws.Range("I122").Value = 1106.8182
$ws.Range("J122").Value = 5610.25
$ws.Range("K122").Value = 3320.4546
$ws.Range("L122").Value = 16830.75
$ws.Range("M122").Value = -870.4546
$ws.Range("N122").Value = -21730.75

# Row 134
$ws.Range("H134").Value = 72856.55499999999
$ws.Range("J134").Value = 72856.55499999999
$ws.Range("L134").Value = 218569.665
$ws.Range("N134").Value = -223639.665

# Row 137
$ws.Range("H137").Value = 46861.43
$ws.Range("J137").Value = 46861.43
$ws.Range("L137").Value = 46861.43
$ws.Range("N137").Value = -57061.43
